# The document contains malformed M2Doc field tokens such as "{m" and
# "{m:" written inside a single Word run (e.g. <w:r><w:t>{m</w:t></w:r>).
# The new TokenIteratorFieldRewriterSplit parser requires the opening
# "{" to live in its own run, separate from the rest of the token text,
# so here we split those runs into "{" + remainder while leaving their
# character formatting untouched.
#
# Word naturally splits a run in two wherever a range boundary is
# introduced inside it (e.g. adding/removing a bookmark on a
# sub-range). We (ab)use that to isolate the "{" character into its
# own run without altering any formatting.

$d = $word.ActiveDocument

$searchStart = 0
$occurrence = 0
$bookmarkIndex = 0

while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("{", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $occurrence = $occurrence + 1
    $bracePos = $rng.Start

    # Occurrence 1 is "{m:for v | self.eClassifiers}" and is left as-is.
    # Occurrence 2 is "{m:v.name}" -> split into "{" + "m:v.name}"...
    # Occurrence 3 is "{m:enfor}"  -> split into "{" + "m:enfor}"...
    if ($occurrence -gt 1) {
        $bookmarkIndex = $bookmarkIndex + 1
        $brace = $d.Range($bracePos, $bracePos + 1)
        $bookmarkName = "m2docSplit" + $bookmarkIndex
        $d.Bookmarks.Add($bookmarkName, $brace)
        $d.Bookmarks($bookmarkName).Delete()
    }

    $searchStart = $bracePos + 1
}
